# Insert a new "ID" column before column B, shifting the existing
# Domain/Dimension/Indicator columns one column to the right, then
# populate the new column with numeric IDs for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns B:D -> C:E by inserting a new blank column at B.
$ws.Columns("B:B").Insert()

# Header for the new ID column (same bold style as the other headers in row 1).
$ws.Range("B1").Value = "ID"
$ws.Range("B1").Font.Bold = $true

# IDs corresponding to each data row (rows with a value in column A).
$ids = @{
    4  = 1001
    5  = 1002
    6  = 1003
    8  = 1004
    9  = 1005
    11 = 1006
    12 = 1007
    13 = 1008
    17 = 2001
    19 = 2002
    20 = 2003
    21 = 2004
    22 = 2005
    26 = 3001
    27 = 3002
    28 = 3003
    29 = 3004
    30 = 3005
    31 = 3006
    32 = 3007
    34 = 3008
    35 = 3009
    39 = 4001
    40 = 4002
    41 = 4003
    42 = 4004
    43 = 4005
    44 = 4006
    46 = 4007
    47 = 4008
    48 = 4009
    50 = 4010
    51 = 4011
    52 = 4012
}

foreach ($row in $ids.Keys) {
    $ws.Cells.Item($row, 2).Value = $ids[$row]
}

$ws.Range("H16").Select()
